$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are plain numeric-looking strings in the source data (e.g. "1.00",
# "505.30"), so a leading apostrophe is used to force Excel to store them as text
# instead of silently converting them to numbers (which would drop formatting such
# as trailing zeros, e.g. "1.00" -> 1).

$ws.Range("D2").Value = "'59.658.63"
$ws.Range("E2").Value = "  +8.19%  "
$ws.Range("D3").Value = "'2.578.61"
$ws.Range("E3").Value = "  +10.02%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'505.30"
$ws.Range("E5").Value = "  +6.25%  "
$ws.Range("D6").Value = "'156.75"
$ws.Range("E6").Value = "  +7.57%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.609"
$ws.Range("E8").Value = "  -4.50%  "
$ws.Range("D9").Value = "'2.577.29"
$ws.Range("E9").Value = "  +9.76%  "
$ws.Range("D10").Value = "'6.10"
$ws.Range("E10").Value = "  +12.12%  "
$ws.Range("E11").Value = "  +6.63%  "
$ws.Range("E12").Value = "  +4.97%  "
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "'3.019.97"
$ws.Range("E14").Value = "  +9.85%  "
$ws.Range("D15").Value = "'59.488.32"
$ws.Range("E15").Value = "  +7.94%  "
$ws.Range("E16").Value = "  +9.26%  "
$ws.Range("E17").Value = "  +5.74%  "
$ws.Range("D18").Value = "'2.573.40"
$ws.Range("E18").Value = "  +9.92%  "
$ws.Range("E19").Value = "  +4.02%  "
$ws.Range("D20").Value = "'339.32"
$ws.Range("E20").Value = "  +7.64%  "
$ws.Range("E21").Value = "  +7.70%  "
$ws.Range("E22").Value = "  +7.39%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +5.54%  "
$ws.Range("E25").Value = "  +5.76%  "
$ws.Range("E26").Value = "  +8.87%  "
$ws.Range("D27").Value = "'2.673.52"
$ws.Range("E27").Value = "  +9.45%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "'0.0₃0832"
$ws.Range("E29").Value = "  +11.66%  "
$ws.Range("D30").Value = "'7.33"
$ws.Range("E30").Value = "  +4.41%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'156.96"
$ws.Range("E32").Value = "  +8.91%  "
$ws.Range("D33").Value = "'19.42"
$ws.Range("E33").Value = "  +6.96%  "
$ws.Range("E34").Value = "  +6.79%  "
$ws.Range("D35").Value = "'5.51"
$ws.Range("E35").Value = "  +7.61%  "
$ws.Range("D36").Value = "'1.19"
$ws.Range("E36").Value = "  +9.89%  "
$ws.Range("D37").Value = "'3.90"
$ws.Range("E37").Value = "  +9.36%  "
$ws.Range("D38").Value = "'0.855"
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("D39").Value = "'301.09"
$ws.Range("E39").Value = "  +20.20%  "
$ws.Range("D40").Value = "'3.72"
$ws.Range("E40").Value = "  +9.10%  "
$ws.Range("E41").Value = "  +8.84%  "
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("D43").Value = "'0.0574"
$ws.Range("E43").Value = "  +10.78%  "
$ws.Range("E44").Value = "  +9.61%  "
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "'0.783"
$ws.Range("E46").Value = "  +25.32%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "'4.91"
$ws.Range("E48").Value = "  +13.46%  "
$ws.Range("D49").Value = "'19.21"
$ws.Range("E49").Value = "  +15.46%  "
$ws.Range("E50").Value = "  +7.59%  "
$ws.Range("D51").Value = "'10.27"
$ws.Range("E51").Value = "  +1.06%  "
